# "error solve ifrs list"
# The IFRS figures in the 현대오토에버 sheet were re-scaled (source data switched
# from KRW to a different reporting unit) and a few now-unused metric columns
# (자산총계/J, 자본금/O, BPS/AD, 현금배당성향/AH) are cleared for the affected
# rows. Apply the corrected figures row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12)
$ws.Range("D2").Value = 11219
$ws.Range("E2").Value = 652
$ws.Range("F2").Value = 652
$ws.Range("G2").Value = 715
$ws.Range("H2").Value = 527
$ws.Range("I2").Value = 524
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 6250
$ws.Range("L2").Value = 3477
$ws.Range("M2").Value = 2772
$ws.Range("N2").Value = 2767
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 103
$ws.Range("Q2").Value = 909
$ws.Range("R2").Value = -126
$ws.Range("S2").Value = -88
$ws.Range("T2").Value = 82
$ws.Range("U2").Value = 828
$ws.Range("V2").Value = 1
$ws.Range("W2").Value = 5.81
$ws.Range("X2").Value = 4.69
$ws.Range("Y2").Value = 21.02
$ws.Range("Z2").Value = 9.5
$ws.Range("AA2").Value = 125.44
$ws.Range("AB2").Value = 2611.84
$ws.Range("AC2").Value = 2609
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").Value = 13399
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = 450
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").Value = 17.72
$ws.Range("AJ2").Value = 20652420

# Row 3 (2015/12)
$ws.Range("D3").Value = 12980
$ws.Range("E3").Value = 682
$ws.Range("F3").Value = 682
$ws.Range("G3").Value = 751
$ws.Range("H3").Value = 536
$ws.Range("I3").Value = 534
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 6873
$ws.Range("L3").Value = 3698
$ws.Range("M3").Value = 3176
$ws.Range("N3").Value = 3172
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 103
$ws.Range("Q3").Value = 111
$ws.Range("R3").Value = -135
$ws.Range("S3").Value = -96
$ws.Range("T3").Value = 93
$ws.Range("U3").Value = 18
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 5.26
$ws.Range("X3").Value = 4.13
$ws.Range("Y3").Value = 18
$ws.Range("Z3").Value = 8.16
$ws.Range("AA3").Value = 116.44
$ws.Range("AB3").Value = 2999.35
$ws.Range("AC3").Value = 2588
$ws.Range("AD3").ClearContents()
$ws.Range("AE3").Value = 15359
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 690
$ws.Range("AH3").ClearContents()
$ws.Range("AI3").Value = 26.66
$ws.Range("AJ3").Value = 20652420

# Row 4 (2016/12)
$ws.Range("D4").Value = 13360
$ws.Range("E4").Value = 804
$ws.Range("F4").Value = 804
$ws.Range("G4").Value = 875
$ws.Range("H4").Value = 647
$ws.Range("I4").Value = 644
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 7961
$ws.Range("L4").Value = 4224
$ws.Range("M4").Value = 3736
$ws.Range("N4").Value = 3709
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 103
$ws.Range("Q4").Value = 686
$ws.Range("R4").Value = -522
$ws.Range("S4").Value = -145
$ws.Range("T4").Value = 294
$ws.Range("U4").Value = 392
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 6.02
$ws.Range("X4").Value = 4.84
$ws.Range("Y4").Value = 18.73
$ws.Range("Z4").Value = 8.720000000000001
$ws.Range("AA4").Value = 113.05
$ws.Range("AB4").Value = 3510.89
$ws.Range("AC4").Value = 3120
$ws.Range("AD4").ClearContents()
$ws.Range("AE4").Value = 17959
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 690
$ws.Range("AH4").ClearContents()
$ws.Range("AI4").Value = 22.12
$ws.Range("AJ4").Value = 20652420

# Row 5 (2017/12)
$ws.Range("D5").Value = 14734
$ws.Range("E5").Value = 729
$ws.Range("F5").Value = 729
$ws.Range("G5").Value = 767
$ws.Range("H5").Value = 552
$ws.Range("I5").Value = 537
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 8040
$ws.Range("L5").Value = 3882
$ws.Range("M5").Value = 4158
$ws.Range("N5").Value = 4115
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 103
$ws.Range("Q5").Value = 500
$ws.Range("R5").Value = 31
$ws.Range("S5").Value = -142
$ws.Range("T5").Value = 145
$ws.Range("U5").Value = 354
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 4.95
$ws.Range("X5").Value = 3.74
$ws.Range("Y5").Value = 13.72
$ws.Range("Z5").Value = 6.9
$ws.Range("AA5").Value = 93.37
$ws.Range("AB5").Value = 3941.58
$ws.Range("AC5").Value = 2600
$ws.Range("AD5").ClearContents()
$ws.Range("AE5").Value = 19927
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 690
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 26.54
$ws.Range("AJ5").Value = 20652420

# Row 6 (2018/12)
$ws.Range("D6").Value = 14249
$ws.Range("E6").Value = 702
$ws.Range("F6").Value = 702
$ws.Range("G6").Value = 751
$ws.Range("H6").Value = 552
$ws.Range("I6").Value = 540
$ws.Range("K6").Value = 8291
$ws.Range("L6").Value = 3775
$ws.Range("M6").Value = 4516
$ws.Range("N6").Value = 4461
$ws.Range("P6").Value = 103
$ws.Range("Q6").Value = 456
$ws.Range("R6").Value = -266
$ws.Range("S6").Value = -93
$ws.Range("T6").Value = 208
$ws.Range("U6").Value = 248
$ws.Range("V6").Value = 51
$ws.Range("W6").Value = 4.92
$ws.Range("X6").Value = 3.88
$ws.Range("Y6").Value = 12.59
$ws.Range("Z6").Value = 6.76
$ws.Range("AA6").Value = 83.59
$ws.Range("AB6").Value = 4269.97
$ws.Range("AC6").Value = 2614
$ws.Range("AD6").ClearContents()
$ws.Range("AE6").Value = 21602
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 690
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 26.4
$ws.Range("AJ6").Value = 20652420

# Row 7 (2019/12(E))
$ws.Range("D7").Value = 15746
$ws.Range("E7").Value = 802
$ws.Range("G7").Value = 855
$ws.Range("H7").Value = 612
$ws.Range("I7").Value = 596
$ws.Range("K7").Value = 9688
$ws.Range("L7").Value = 4620
$ws.Range("M7").Value = 5068
$ws.Range("N7").Value = 5004
$ws.Range("P7").Value = 102
$ws.Range("Q7").Value = 832
$ws.Range("R7").Value = -430
$ws.Range("S7").Value = -132
$ws.Range("T7").Value = 358
$ws.Range("U7").Value = 360
$ws.Range("W7").Value = 5.09
$ws.Range("X7").Value = 3.89
$ws.Range("Y7").Value = 12.61
$ws.Range("Z7").Value = 6.81
$ws.Range("AA7").Value = 91.17
$ws.Range("AC7").Value = 2851
$ws.Range("AD7").Value = 17.19
$ws.Range("AE7").Value = 23826
$ws.Range("AF7").Value = 2.06
$ws.Range("AG7").Value = 715
$ws.Range("AH7").Value = 1.46
$ws.Range("AI7").Value = 25.17

# Row 8 (2020/12(E))
$ws.Range("D8").Value = 17738
$ws.Range("E8").Value = 939
$ws.Range("G8").Value = 1004
$ws.Range("H8").Value = 720
$ws.Range("I8").Value = 702
$ws.Range("K8").Value = 10684
$ws.Range("L8").Value = 5058
$ws.Range("M8").Value = 5626
$ws.Range("N8").Value = 5548
$ws.Range("P8").Value = 102
$ws.Range("Q8").Value = 1054
$ws.Range("R8").Value = -532
$ws.Range("S8").Value = -152
$ws.Range("T8").Value = 500
$ws.Range("U8").Value = 560
$ws.Range("W8").Value = 5.29
$ws.Range("X8").Value = 4.06
$ws.Range("Y8").Value = 13.31
$ws.Range("Z8").Value = 7.07
$ws.Range("AA8").Value = 89.91
$ws.Range("AC8").Value = 3343
$ws.Range("AD8").Value = 13.7
$ws.Range("AE8").Value = 26421
$ws.Range("AF8").Value = 1.73
$ws.Range("AG8").Value = 765
$ws.Range("AH8").Value = 1.67
$ws.Range("AI8").Value = 22.88

# Row 9 (2021/12(E))
$ws.Range("D9").Value = 19918
$ws.Range("E9").Value = 1071
$ws.Range("G9").Value = 1152
$ws.Range("H9").Value = 826
$ws.Range("I9").Value = 808
$ws.Range("K9").Value = 11890
$ws.Range("L9").Value = 5610
$ws.Range("M9").Value = 6280
$ws.Range("N9").Value = 6183
$ws.Range("P9").Value = 102
$ws.Range("Q9").Value = 1162
$ws.Range("R9").Value = -534
$ws.Range("S9").Value = -162
$ws.Range("T9").Value = 500
$ws.Range("U9").Value = 650
$ws.Range("W9").Value = 5.38
$ws.Range("X9").Value = 4.15
$ws.Range("Y9").Value = 13.77
$ws.Range("Z9").Value = 7.32
$ws.Range("AA9").Value = 89.34
$ws.Range("AC9").Value = 3845
$ws.Range("AD9").Value = 11.91
$ws.Range("AE9").Value = 29443
$ws.Range("AF9").Value = 1.56
$ws.Range("AG9").Value = 815
$ws.Range("AH9").Value = 1.78
$ws.Range("AI9").Value = 21.2
